$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.637.16"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "1.842.35"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'259.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("D8").Value = "'0.3149"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.43%  "
$ws.Range("D9").Value = "'0.06801"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").Value = "'18.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "'0.7809"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "'0.07756"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "1.833.80"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "'87.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'0.9997"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'13.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.9997"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "26.657.84"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "2.075.80"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "'4.606"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "'5.973"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "'9.322"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.22%  "
$ws.Range("D25").Value = "'142.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").Value = "'2.210"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").Value = "'1.681"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("D28").Value = "'16.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "'110.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").Value = "'4.188"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "'0.08727"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").Value = "'4.073"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").Value = "'0.7319"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").Value = "'1.140"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("D36").Value = "'2.859"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Value = "'3.091"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'2.261"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("D39").Value = "'0.01732"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("D40").Value = "'0.4802"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").Value = "'0.8944"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("D42").Value = "'109.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").Value = "'0.9999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'7.672"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").Value = "'0.4163"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "'8.995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "'0.1238"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").Value = "'0.05818"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("D50").Value = "'34.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").Value = "'0.8932"
$ws.Range("D51").Style = "Normal"
